$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.897.80"
$ws.Range("E2").Value = "  +1.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.411.57"
$ws.Range("E3").Value = "  +1.63%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.19"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.68"
$ws.Range("E6").Value = "  +2.62%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.417.90"
$ws.Range("E8").Value = "  +1.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.554"
$ws.Range("E9").Value = "  +5.02%  "

$ws.Range("E10").Value = "  -1.40%  "

$ws.Range("E11").Value = "  +4.47%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.434"
$ws.Range("E12").Value = "  -1.06%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.000.63"
$ws.Range("E13").Value = "  +1.73%  "

$ws.Range("E14").Value = "  -3.34%  "

$ws.Range("E15").Value = "  +6.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.47"
$ws.Range("E16").Value = "  +2.40%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.871.83"
$ws.Range("E17").Value = "  +1.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.400.05"
$ws.Range("E18").Value = "  +1.38%  "

$ws.Range("E19").Value = "  -0.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.10"
$ws.Range("E20").Value = "  +2.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "382.39"
$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.07"
$ws.Range("E22").Value = "  -3.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.16"
$ws.Range("E24").Value = "  +2.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.533"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("E26").Value = "  +24.44%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.49"
$ws.Range("E27").Value = "  +2.02%  "

$ws.Range("E28").Value = "  -0.12%  "

$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("E30").Value = "  +11.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.38"
$ws.Range("E31").Value = "  +5.92%  "

$ws.Range("E32").Value = "  -0.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.48"
$ws.Range("E33").Value = "  +2.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.30"
$ws.Range("E34").Value = "  +1.32%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.83"
$ws.Range("E36").Value = "  +1.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.17"
$ws.Range("E37").Value = "  +0.84%  "

$ws.Range("E38").Value = "  +0.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.950.58"
$ws.Range("E39").Value = "  +6.24%  "

$ws.Range("E40").Value = "  +2.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.83"
$ws.Range("E41").Value = "  -1.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.77"
$ws.Range("E42").Value = "  -0.81%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0314"
$ws.Range("E43").Value = "  -4.50%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.08"
$ws.Range("E44").Value = "  +2.75%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.34"
$ws.Range("E45").Value = "  +1.66%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.759"
$ws.Range("E46").Value = "  +1.89%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.29"
$ws.Range("E47").Value = "  +6.40%  "

$ws.Range("E48").Value = "  +2.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.21"
$ws.Range("E49").Value = "  +21.77%  "

$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.837"
$ws.Range("E50").Value = "  +3.75%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.38"
$ws.Range("E51").Value = "  +0.93%  "
